$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataBF = New-Object "object[,]" 24,5
$dataBF[0,0] = 0.891281037216288
$dataBF[0,1] = 0.227868235862644
$dataBF[0,2] = 0.008829382088016757
$dataBF[0,3] = 0.424636301865192
$dataBF[0,4] = 0.511900057648603
$dataBF[1,0] = 0.7821342540257206
$dataBF[1,1] = 0.2013612742698854
$dataBF[1,2] = 0.00796129246657884
$dataBF[1,3] = 0.3703224242887444
$dataBF[1,4] = 0.4991551764899711
$dataBF[2,0] = 0.715054343910424
$dataBF[2,1] = 0.1850036139169049
$dataBF[2,2] = 0.007425584329315171
$dataBF[2,3] = 0.3370727147557346
$dataBF[2,4] = 0.4917844073917266
$dataBF[3,0] = 0.6877035934890046
$dataBF[3,1] = 0.1783172446287438
$dataBF[3,2] = 0.007206620249466056
$dataBF[3,3] = 0.3235461137417133
$dataBF[3,4] = 0.4888945045365602
$dataBF[4,0] = 0.683161131083466
$dataBF[4,1] = 0.1772057465150567
$dataBF[4,2] = 0.007170222232623757
$dataBF[4,3] = 0.3213013576355337
$dataBF[4,4] = 0.488421491470362
$dataBF[5,0] = 0.7146855423868601
$dataBF[5,1] = 0.1849135218971583
$dataBF[5,2] = 0.007422633943278356
$dataBF[5,3] = 0.3368902001133449
$dataBF[5,4] = 0.4917449733671404
$dataBF[6,0] = 0.8536608088396633
$dataBF[6,1] = 0.2187458238674935
$dataBF[6,2] = 0.008530636135073877
$dataBF[6,3] = 0.4058870429058317
$dataBF[6,4] = 0.5074109087968353
$dataBF[7,0] = 1.125667569260031
$dataBF[7,1] = 0.2844324375343206
$dataBF[7,2] = 0.01068126729409613
$dataBF[7,3] = 0.5420748521896286
$dataBF[7,4] = 0.5417675613070116
$dataBF[8,0] = 1.325184521455583
$dataBF[8,1] = 0.3322882203825657
$dataBF[8,2] = 0.01224697816459752
$dataBF[8,3] = 0.642816469174079
$dataBF[8,4] = 0.5692687019538027
$dataBF[9,0] = 1.415878935840453
$dataBF[9,1] = 0.3539710371430829
$dataBF[9,2] = 0.01295597692462991
$dataBF[9,3] = 0.6888270794604097
$dataBF[9,4] = 0.5822797367671626
$dataBF[10,0] = 1.450212479127003
$dataBF[10,1] = 0.3621691174283228
$dataBF[10,2] = 0.01322397201867176
$dataBF[10,3] = 0.7062788127806243
$dataBF[10,4] = 0.5872793665215852
$dataBF[11,0] = 1.442818615038959
$dataBF[11,1] = 0.3604040838560536
$dataBF[11,2] = 0.01316627646838953
$dataBF[11,3] = 0.7025189657318691
$dataBF[11,4] = 0.5861993661312965
$dataBF[12,0] = 1.418703791635664
$dataBF[12,1] = 0.3546457545478177
$dataBF[12,2] = 0.01297803492241911
$dataBF[12,3] = 0.6902622597063157
$dataBF[12,4] = 0.5826895999395845
$dataBF[13,0] = 1.403931378352922
$dataBF[13,1] = 0.3511169470808397
$dataBF[13,2] = 0.0128626674887542
$dataBF[13,3] = 0.6827584534899103
$dataBF[13,4] = 0.5805492449541987
$dataBF[14,0] = 1.319256010063896
$dataBF[14,1] = 0.3308694275111179
$dataBF[14,2] = 0.01220057630517601
$dataBF[14,3] = 0.639813435357766
$dataBF[14,4] = 0.5684285290814017
$dataBF[15,0] = 1.267292684093718
$dataBF[15,1] = 0.3184257859094828
$dataBF[15,2] = 0.01179355792734071
$dataBF[15,3] = 0.6135164451768702
$dataBF[15,4] = 0.5611215421691469
$dataBF[16,0] = 1.23739851142949
$dataBF[16,1] = 0.3112603546215382
$dataBF[16,2] = 0.01155914689276472
$dataBF[16,3] = 0.5984081774049912
$dataBF[16,4] = 0.5569658393738024
$dataBF[17,0] = 1.227275812269681
$dataBF[17,1] = 0.3088328648848915
$dataBF[17,2] = 0.01147972760001181
$dataBF[17,3] = 0.5932956408562973
$dataBF[17,4] = 0.5555668567712644
$dataBF[18,0] = 1.272824921058884
$dataBF[18,1] = 0.3197512807541898
$dataBF[18,2] = 0.01183691742154025
$dataBF[18,3] = 0.6163140223165016
$dataBF[18,4] = 0.5618945058666469
$dataBF[19,0] = 1.425787190507151
$dataBF[19,1] = 0.3563374627162545
$dataBF[19,2] = 0.01303333941032747
$dataBF[19,3] = 0.6938615595734632
$dataBF[19,4] = 0.5837185269568863
$dataBF[20,0] = 1.525695502703115
$dataBF[20,1] = 0.3801743374466469
$dataBF[20,2] = 0.01381242070084454
$dataBF[20,3] = 0.7447106142914066
$dataBF[20,4] = 0.598405449573761
$dataBF[21,0] = 1.472378449820724
$dataBF[21,1] = 0.3674590195964242
$dataBF[21,2] = 0.01339687757999286
$dataBF[21,3] = 0.7175554673538045
$dataBF[21,4] = 0.5905277838387093
$dataBF[22,0] = 1.270323858226448
$dataBF[22,1] = 0.3191520600789204
$dataBF[22,2] = 0.01181731587243462
$dataBF[22,3] = 0.6150492060403536
$dataBF[22,4] = 0.5615449082929018
$dataBF[23,0] = 1.052139364296011
$dataBF[23,1] = 0.2667332554221389
$dataBF[23,2] = 0.01010192572474011
$dataBF[23,3] = 0.5051227919026502
$dataBF[23,4] = 0.5320795236704896
$ws.Range("B2:F25").Value = $dataBF

$dataNO = New-Object "object[,]" 24,2
$dataNO[0,0] = 0.8201990362995204
$dataNO[0,1] = 1.626294133104693
$dataNO[1,0] = 0.8269554030532333
$dataNO[1,1] = 1.597674074346656
$dataNO[2,0] = 0.8315100509956181
$dataNO[2,1] = 1.581611426959796
$dataNO[3,0] = 0.8334682874417609
$dataNO[3,1] = 1.575443445801824
$dataNO[4,0] = 0.8337996237252412
$dataNO[4,1] = 1.574442007655904
$dataNO[5,0] = 0.8315360466860326
$dataNO[5,1] = 1.581526716955437
$dataNO[6,0] = 0.822444341291046
$dataNO[6,1] = 1.616111236584089
$dataNO[7,0] = 0.8078376126043167
$dataNO[7,1] = 1.696013581328145
$dataNO[8,0] = 0.7990697030542435
$dataNO[8,1] = 1.762230245765721
$dataNO[9,0] = 0.7955073947389835
$dataNO[9,1] = 1.79401746244065
$dataNO[10,0] = 0.7942197524327455
$dataNO[10,1] = 1.806296323695108
$dataNO[11,0] = 0.7944943411547314
$dataNO[11,1] = 1.803641063067687
$dataNO[12,0] = 0.7954002303527545
$dataNO[12,1] = 1.795022793749439
$dataNO[13,0] = 0.7959631012151078
$dataNO[13,1] = 1.789775404628529
$dataNO[14,0] = 0.7993110881807439
$dataNO[14,1] = 1.760186576791938
$dataNO[15,0] = 0.8014741645121291
$dataNO[15,1] = 1.74246286157171
$dataNO[16,0] = 0.8027584278971176
$dataNO[16,1] = 1.732425220543718
$dataNO[17,0] = 0.8032001470647003
$dataNO[17,1] = 1.729053468574108
$dataNO[18,0] = 0.8012397489807483
$dataNO[18,1] = 1.744333360078087
$dataNO[19,0] = 0.7951324840674232
$dataNO[19,1] = 1.797547609276052
$dataNO[20,0] = 0.7914985108803947
$dataNO[20,1] = 1.833736282588887
$dataNO[21,0] = 0.7934053094209403
$dataNO[21,1] = 1.814291914354243
$dataNO[22,0] = 0.8013456015225415
$dataNO[22,1] = 1.743487234541078
$dataNO[23,0] = 0.8114442649019864
$dataNO[23,1] = 1.673089269594783
$ws.Range("N2:O25").Value = $dataNO

Write-Host "Applied case with 380 kV values"